$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.775.61'
$ws.Range("E2").Value = '  +2.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.813.78'
$ws.Range("E3").Value = '  +1.25%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '663.16'
$ws.Range("E5").Value = '  +6.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.33'
$ws.Range("E6").Value = '  +3.31%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.810.86'
$ws.Range("E7").Value = '  +1.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  +1.62%  '

$ws.Range("E10").Value = '  +1.00%  '

$ws.Range("E11").Value = '  +2.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.95'
$ws.Range("E12").Value = '  +5.08%  '

$ws.Range("E13").Value = '  -0.86%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.69'
$ws.Range("E14").Value = '  +1.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.456.77'
$ws.Range("E15").Value = '  +1.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.818.66'
$ws.Range("E16").Value = '  +0.15%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.747.63'
$ws.Range("E17").Value = '  +2.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.79'
$ws.Range("E18").Value = '  +0.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.15'
$ws.Range("E19").Value = '  +1.14%  '

$ws.Range("E20").Value = '  +1.00%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '478.30'
$ws.Range("E21").Value = '  +2.33%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.33'
$ws.Range("E22").Value = '  +7.94%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.713'
$ws.Range("E23").Value = '  +1.90%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000146'
$ws.Range("E24").Value = '  -1.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.82'
$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.25'
$ws.Range("E26").Value = '  +2.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.38'
$ws.Range("E27").Value = '  +4.27%  '

$ws.Range("E28").Value = '  -1.01%  '

$ws.Range("E29").Value = '  +0.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.964.55'
$ws.Range("E30").Value = '  +1.27%  '

$ws.Range("E31").Value = '  +7.38%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.32'
$ws.Range("E32").Value = '  +3.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.44'
$ws.Range("E33").Value = '  +2.48%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.43'
$ws.Range("E34").Value = '  +2.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.179'
$ws.Range("E35").Value = '  +16.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.769.19'
$ws.Range("E36").Value = '  +1.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.09'
$ws.Range("E38").Value = '  +1.60%  '

$ws.Range("E39").Value = '  -0.26%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.41'
$ws.Range("E40").Value = '  +2.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.93'
$ws.Range("E41").Value = '  +2.83%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.969'
$ws.Range("E42").Value = '  +0.50%  '

$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.10'
$ws.Range("E44").Value = '  +10.02%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.40'
$ws.Range("E46").Value = '  +5.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '158.84'
$ws.Range("E47").Value = '  +3.35%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '47.80'
$ws.Range("E48").Value = '  +2.44%  '

$ws.Range("E49").Value = '  +0.96%  '

$ws.Range("E50").Value = '  +4.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.50'
$ws.Range("E51").Value = '  +1.33%  '
